$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C10").Value = 0.5499551136970593
$ws.Range("E10").Value = 0.6576019635478055
$ws.Range("G10").Value = 7.735736368705187
$ws.Range("H10").Value = 39.52096861246186
$ws.Range("I10").Value = 52.74329501883294
$ws.Range("J10").Value = 47.25670498116705
$ws.Range("K10").Value = 52.74329501883294
$ws.Range("B11").Value = 1.614989422475507
$ws.Range("E11").Value = 3.927972263761085
$ws.Range("G11").Value = 30.90647844926512
$ws.Range("H11").Value = 44.26416256531046
$ws.Range("I11").Value = 24.82935898542442
$ws.Range("J11").Value = 75.17064101457558
$ws.Range("K11").Value = 24.82935898542442
$ws.Range("C12").Value = 3.80766260903579
$ws.Range("E12").Value = 4.826385040842168
$ws.Range("G12").Value = 17.57130444472443
$ws.Range("H12").Value = 65.67598477980499
$ws.Range("I12").Value = 16.75271077547058
$ws.Range("J12").Value = 83.24728922452942
$ws.Range("K12").Value = 16.75271077547058
$ws.Range("D15").Value = 4.588987080161871
$ws.Range("F15").Value = 4.588987080161871
$ws.Range("G15").Value = 38.90185099316687
$ws.Range("H15").Value = 21.06747471048451
$ws.Range("I15").Value = 40.03067429634861
$ws.Range("J15").Value = 59.96932570365138
$ws.Range("K15").Value = 40.03067429634861
$ws.Range("B16").Value = 0.8860581650638293
$ws.Range("E16").Value = 1.460640766159236
$ws.Range("G16").Value = 51.87160332039483
$ws.Range("H16").Value = 33.63720569819994
$ws.Range("I16").Value = 14.49119098140523
$ws.Range("J16").Value = 85.50880901859477
$ws.Range("K16").Value = 14.49119098140523
$ws.Range("C22").Value = 0.02118677
$ws.Range("D22").Value = 0.00776831
$ws.Range("E22").Value = 0.04865666
$ws.Range("F22").Value = 0.00776831
$ws.Range("G22").Value = 48.68392486517937
$ws.Range("H22").Value = 37.54857113791996
$ws.Range("I22").Value = 13.76750399690066
$ws.Range("J22").Value = 86.23249600309933
$ws.Range("K22").Value = 13.76750399690066
$ws.Range("B23").Value = 18.50934330204853
$ws.Range("C23").Value = 33.69065111221761
$ws.Range("D23").Value = 17.59157887518554
$ws.Range("E23").Value = 52.19999441426614
$ws.Range("F23").Value = 17.59157887518554
$ws.Range("G23").Value = 26.52088558783936
$ws.Range("H23").Value = 48.27323632967825
$ws.Range("I23").Value = 25.2058780824824
$ws.Range("J23").Value = 74.7941219175176
$ws.Range("K23").Value = 25.2058780824824
